# "Simplify report generation window"
#
# The monthly report sheet used to reserve 10 columns (A:J) for a
# Фамилия/Имя/.../Стоимость.../Общая стоимость table. The cost columns
# (H:J - "Стоимость обедов", "Стоимость ужинов", "Общая стоимость") are
# dropped from the generated report header, so the window only needs to
# show columns A:G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The title row (row 1) still spans columns A:J as far as the sheet
# dimension is concerned, but H1:J1 are no longer part of the merged /
# centered title banner - they keep the same bold/italic title font,
# just without the centered horizontal alignment (vertical centering
# only).
$ws.Range("H1:J1").HorizontalAlignment = 1

# The three now-unused column headers (Стоимость обедов / Стоимость
# ужинов / Общая стоимость) are removed entirely from the header row.
$ws.Range("H2:J2").Clear()

# Shrink the title merge from A1:J1 down to A1:G1 to match the smaller
# table, and select the new merged title range.
$ws.Range("A1:J1").UnMerge()
$ws.Range("A1:G1").Merge()
$ws.Range("A1:G1").Select()
